# Update the "Date" metadata value on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-03-24T15:01:17-04:00"

# Add the new "Include from LOINC" worksheet at the end of the workbook.
$mmComm = $wb.Worksheets.Item("Include from Multi-Modal Comm")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Include from LOINC"

# Column widths matching the "Include from Multi-Modal Comm" sheet
# (30.703125 / 50.703125 "raw" OOXML width units; the closest values the
# Excel column-width-in-characters model can reproduce are used here).
$ws.Columns.Item(1).ColumnWidth = 29.75
$ws.Columns.Item(2).ColumnWidth = 49.75

# Header row.
$ws.Range("A1").Value = "Concept"
$ws.Range("B1").Value = "Description"

# LOINC concept rows.
$ws.Range("A2").Value = "99829-4"
$ws.Range("B2").Value = "How often does the individual convey simple messages that are meaningful related to routine daily activities in LOW demand situations?"

$ws.Range("A3").Value = "99830-2"
$ws.Range("B3").Value = "How often does the individual participate in short structured conversations that are meaningful in LOW demand situations?"

$ws.Range("A4").Value = "99831-0"
$ws.Range("B4").Value = "How often does the individual convey complex messages that are meaningful in LOW demand situations?"

$ws.Range("A5").Value = "99832-8"
$ws.Range("B5").Value = "How often does the individual convey simple messages that are meaningful related to routine daily activities in HIGH demand situations?"

$ws.Range("A6").Value = "99833-6"
$ws.Range("B6").Value = "How often does the individual participate in short structured conversations that are meaningful in HIGH demand situations?"

$ws.Range("A7").Value = "99834-4"
$ws.Range("B7").Value = "How often does the individual convey complex messages that are meaningful in HIGH demand situations?"

$ws.Range("A8").Value = "99835-1"
$ws.Range("B8").Value = "Functional Communication Measure - Multi-Modal Functional Communication score [ASHA NOMS]"

# Blank spacer row.
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = ""

# System URI row.
$ws.Range("A10").Value = "System URI"
$ws.Range("B10").Value = "http://loinc.org"

# Copy cell formatting from existing sheets so the new sheet's styling
# matches (bold filled header row, bordered body cells with top-aligned
# wrapped text). The "Metadata" sheet's A1:B1 is used as the header-format
# source since both of its columns carry the header style; "Include from
# Multi-Modal Comm"'s A4:B4 is used for the body-format source since it is
# the first row on that sheet with both columns populated (row 1 there only
# has column A set).
$meta.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$mmComm.Range("A4:B4").Copy()
$ws.Range("A2:B10").PasteSpecial(-4122)

$excel.CutCopyMode = $false
